$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels: "height" and "weight" are inserted (columns E and F),
# and the former "fantasy points" column is relocated to column G.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"
$ws.Range("G1").Value = "fantasy points"

# Match the header formatting already used by B1:E1 (bold, bordered,
# centered) on the two newly added header cells.
$ws.Range("D1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Capture the existing "fantasy points" values (currently in column E) before
# overwriting column E with the new height data.
$fantasyPoints = @()
for ($r = 2; $r -le 17; $r++) {
    $fantasyPoints += $ws.Cells.Item($r, 5).Value2
}

# Populate height (E) and weight (F) for every data row, then move the
# captured fantasy-points values into the new column G.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 253
    $ws.Cells.Item($r, 7).Value = $fantasyPoints[$r - 2]
}
